$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows to fill with 5 in columns C and D (per diff)
$rows = @(8, 10, 11, 12, 13, 14, 15, 16, 17, 19, 20, 22)

foreach ($r in $rows) {
    $ws.Cells.Item($r, 3).Value = 5
    $ws.Cells.Item($r, 4).Value = 5
}

# Update the active selection to C8 (matches final saved selection in the sheet view)
$ws.Range("C8").Select()
